# Weekly update: insert a new price record at row 40, pushing the
# existing rows 40-96 down to 41-97 (matches the source diff exactly).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row before the current row 40, shifting rows 40:96 -> 41:97
$ws.Rows("40:40").Insert()

# Populate the newly inserted row 40 with the new record's data
$ws.Cells(40, 1).Value = 11
$ws.Cells(40, 2).Value = "Vega Monumental Concepción"
$ws.Cells(40, 3).Value = "Bíobío"
$ws.Cells(40, 4).Value = 44665
$ws.Cells(40, 5).Value = 8
$ws.Cells(40, 6).Value = 100112021
$ws.Cells(40, 7).Value = "Ají"
$ws.Cells(40, 8).Value = "Americana (o)"
$ws.Cells(40, 9).Value = "Primera"
$ws.Cells(40, 10).Value = 100
$ws.Cells(40, 11).Value = 28000
$ws.Cells(40, 12).Value = 30000
$ws.Cells(40, 13).Value = 29000
$ws.Cells(40, 14).Value = "`$/caja 25 kilos"
$ws.Cells(40, 15).Value = "Provincia de Limarí"
$ws.Cells(40, 16).Value = 1160
$ws.Cells(40, 17).Value = 25
$ws.Cells(40, 18).Value = "Hortaliza"
